$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1159.8857
$ws.Range("I15").Value = 1159.8857
$ws.Range("K15").Value = 3479.6571
$ws.Range("M15").Value = -3310.6571
$ws.Range("H33").Value = 383.66666
$ws.Range("I33").Value = 369
$ws.Range("J33").Value = 457
$ws.Range("K33").Value = 369
$ws.Range("L33").Value = 457
$ws.Range("M33").Value = -140
$ws.Range("N33").Value = -915
$ws.Range("H99").Value = 4146.8
$ws.Range("J99").Value = 5143.125
$ws.Range("L99").Value = 15429.375
$ws.Range("N99").Value = -18425.375
$ws.Range("H132").Value = 6204.2593
$ws.Range("I132").Value = 3795.2778
$ws.Range("K132").Value = 11385.8334
$ws.Range("M132").Value = -8855.8334
$ws.Range("H137").Value = 2052.9524
$ws.Range("I137").Value = 1945.4375
$ws.Range("K137").Value = 5836.3125
$ws.Range("M137").Value = -3286.3125

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2405.4827
$ws.Range("I2").Value = 2648.45
$ws.Range("J2").Value = 1865.5555
$ws.Range("K2").Value = 2648.45
$ws.Range("L2").Value = 1865.5555
$ws.Range("M2").Value = -2535.45
$ws.Range("N2").Value = -2091.5555
$ws.Range("H19").Value = 35000
$ws.Range("I19").Value = 35000
$ws.Range("K19").Value = 35000
$ws.Range("M19").Value = -34771
$ws.Range("H43").Value = 35570.332
$ws.Range("J43").Value = 34534.715
$ws.Range("L43").Value = 34534.715
$ws.Range("N43").Value = -35160.715
$ws.Range("H97").Value = 1018.7857
$ws.Range("I97").Value = 986.0833
$ws.Range("J97").Value = 1215
$ws.Range("K97").Value = 986.0833
$ws.Range("L97").Value = 1215
$ws.Range("M97").Value = -490.0833
$ws.Range("N97").Value = -2207
$ws.Range("H116").Value = 2405.4827
$ws.Range("I116").Value = 2648.45
$ws.Range("J116").Value = 1865.5555
$ws.Range("K116").Value = 2648.45
$ws.Range("L116").Value = 1865.5555
$ws.Range("M116").Value = -354.4499999999998
$ws.Range("N116").Value = -6453.5555
$ws.Range("H132").Value = 4174342
$ws.Range("I132").Value = 5899
$ws.Range("J132").Value = 16679671
$ws.Range("K132").Value = 17697
$ws.Range("L132").Value = 50039013
$ws.Range("M132").Value = -15167
$ws.Range("N132").Value = -50044073

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2405.4827
$ws.Range("I3").Value = 2648.45
$ws.Range("J3").Value = 1865.5555
$ws.Range("K3").Value = 2648.45
$ws.Range("L3").Value = 1865.5555
$ws.Range("M3").Value = -2534.45
$ws.Range("N3").Value = -2093.5555
$ws.Range("H20").Value = 19479.2
$ws.Range("I20").Value = 19479.2
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 19479.2
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -19232.2
$ws.Range("N20").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4547860
$ws.Range("I16").Value = 5002376
$ws.Range("K16").Value = 5002376
$ws.Range("M16").Value = -5002089
$ws.Range("H58").Value = 2652.8696
$ws.Range("I58").Value = 2656.4
$ws.Range("K58").Value = 2656.4
$ws.Range("M58").Value = -2453.4
$ws.Range("H94").Value = 1144.25
$ws.Range("I94").Value = 1469.3334
$ws.Range("K94").Value = 1469.3334
$ws.Range("M94").Value = -1018.3334
$ws.Range("H105").Value = 1836.7391
$ws.Range("I105").Value = 1512.2
$ws.Range("J105").Value = 2445.25
$ws.Range("K105").Value = 1512.2
$ws.Range("L105").Value = 2445.25
$ws.Range("M105").Value = 234.8
$ws.Range("N105").Value = -5939.25
$ws.Range("H107").Value = 1327.3334
$ws.Range("I107").Value = 875.35
$ws.Range("J107").Value = 3587.25
$ws.Range("K107").Value = 875.35
$ws.Range("L107").Value = 3587.25
$ws.Range("M107").Value = 1044.65
$ws.Range("N107").Value = -7427.25
$ws.Range("H113").Value = 4547860
$ws.Range("I113").Value = 5002376
$ws.Range("K113").Value = 5002376
$ws.Range("M113").Value = -5000206
$ws.Range("H136").Value = 2652.8696
$ws.Range("I136").Value = 2656.4
$ws.Range("K136").Value = 7969.200000000001
$ws.Range("M136").Value = -5419.200000000001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 30547.166
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 30547.166
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 91641.49800000001
$ws.Range("M63").ClearContents()
$ws.Range("N63").Value = -93139.49800000001
$ws.Range("H66").Value = 30547.166
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 30547.166
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 274924.494
$ws.Range("M66").ClearContents()
$ws.Range("N66").Value = -282412.494
$ws.Range("H131").Value = 4015.0476
$ws.Range("I131").Value = 2305.75
$ws.Range("K131").Value = 6917.25
$ws.Range("M131").Value = -1877.25
$ws.Range("H132").Value = 1383.1428
$ws.Range("I132").Value = 855.4
$ws.Range("K132").Value = 7698.599999999999
$ws.Range("M132").Value = -5168.599999999999
$ws.Range("H137").Value = 5196.077
$ws.Range("I137").Value = 1763
$ws.Range("J137").Value = 8629.154
$ws.Range("K137").Value = 5289
$ws.Range("L137").Value = 25887.462
$ws.Range("M137").Value = -189
$ws.Range("N137").Value = -36087.462

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 505.86365
$ws.Range("J97").Value = 156
$ws.Range("L97").Value = 156
$ws.Range("N97").Value = -1148

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 4632008.5
$ws.Range("I68").Value = 8335105
$ws.Range("J68").Value = 3137.75
$ws.Range("K68").Value = 8335105
$ws.Range("L68").Value = 3137.75
$ws.Range("M68").Value = -8334356
$ws.Range("N68").Value = -4635.75
$ws.Range("H71").Value = 4632008.5
$ws.Range("I71").Value = 8335105
$ws.Range("J71").Value = 3137.75
$ws.Range("K71").Value = 41675525
$ws.Range("L71").Value = 15688.75
$ws.Range("M71").Value = -41671781
$ws.Range("N71").Value = -23176.75
$ws.Range("H82").Value = 5581.1333
$ws.Range("J82").Value = 8499.857
$ws.Range("L82").Value = 8499.857
$ws.Range("N82").Value = -9221.857
$ws.Range("H85").Value = 5581.1333
$ws.Range("J85").Value = 8499.857
$ws.Range("L85").Value = 8499.857
$ws.Range("N85").Value = -10995.857
$ws.Range("H93").Value = 2319079.8
$ws.Range("I93").Value = 2205.125
$ws.Range("J93").Value = 6952829
$ws.Range("K93").Value = 2205.125
$ws.Range("L93").Value = 6952829
$ws.Range("M93").Value = -957.125
$ws.Range("N93").Value = -6955325
$ws.Range("H134").Value = 95832
$ws.Range("J134").Value = 95832
$ws.Range("L134").Value = 95832
$ws.Range("N134").Value = -105972
$ws.Range("H136").Value = 3103
$ws.Range("J136").Value = 3844.25
$ws.Range("L136").Value = 11532.75
$ws.Range("N136").Value = -16632.75
$ws.Range("H138").Value = 74999.5
$ws.Range("J138").Value = 74999.5
$ws.Range("L138").Value = 74999.5
$ws.Range("N138").Value = -85279.5
$ws.Range("H140").Value = 47713.145
$ws.Range("J140").Value = 47713.145
$ws.Range("L140").Value = 47713.145
$ws.Range("N140").Value = -58073.145
$ws.Range("H141").Value = 100799.2
$ws.Range("J141").Value = 100799.2
$ws.Range("L141").Value = 100799.2
$ws.Range("N141").Value = -111159.2

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 16993.857
$ws.Range("J41").Value = 15327
$ws.Range("L41").Value = 15327
$ws.Range("N41").Value = -16107
$ws.Range("H75").Value = 49073.668
$ws.Range("I75").Value = 23611
$ws.Range("K75").Value = 23611
$ws.Range("M75").Value = -22675
$ws.Range("H78").Value = 49073.668
$ws.Range("I78").Value = 23611
$ws.Range("K78").Value = 70833
$ws.Range("M78").Value = -66153
$ws.Range("H81").Value = 2246.7
$ws.Range("I81").Value = 1893.8334
$ws.Range("J81").Value = 2776
$ws.Range("K81").Value = 3787.6668
$ws.Range("L81").Value = 5552
$ws.Range("M81").Value = -2726.6668
$ws.Range("N81").Value = -7674
$ws.Range("H84").Value = 2246.7
$ws.Range("I84").Value = 1893.8334
$ws.Range("J84").Value = 2776
$ws.Range("K84").Value = 18938.334
$ws.Range("L84").Value = 27760
$ws.Range("M84").Value = -13634.334
$ws.Range("N84").Value = -38368
$ws.Range("H96").Value = 9776.799999999999
$ws.Range("I96").Value = 8346
$ws.Range("J96").Value = 15500
$ws.Range("K96").Value = 8346
$ws.Range("L96").Value = 15500
$ws.Range("M96").Value = -6973
$ws.Range("N96").Value = -18246
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()
$ws.Range("H133").Value = 60000
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
